# Update "Hoja1" worksheet: refresh the printed date and the three price
# values in the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Date shown at the top of the sheet (A1): 24/04/2024 -> 24/05/2024
$ws.Range("A1").Value = 45436

# Updated prices for the three hinge items (column D, rows 26-28)
$ws.Range("D26").Value = 437.461
$ws.Range("D27").Value = 504.071
$ws.Range("D28").Value = 618.837
